# LaptopManagement: load data from Excel file and SQL server.
# ProductDate column (D) is converted from free-text strings to real Excel
# dates so the app can read/sort them natively; a handful of blank,
# date-formatted placeholder rows are appended below the table for future
# SQL-sourced rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Convert the ProductDate column to real dates -------------------------
# Header keeps its bold/border look but now also carries the date format.
$ws.Range("D1").NumberFormat = "dd/mm/yyyy"

$ws.Range("D2:D6").NumberFormat = "dd/mm/yyyy"
$ws.Range("D2").Value = 45000   # 2023-03-15
$ws.Range("D3").Value = 44155
$ws.Range("D4").Value = 44931   # 2023-01-05
$ws.Range("D5").Value = 44822   # 2022-09-18
$ws.Range("D6").Value = 45117   # 2023-07-10

# --- Reserve a few more date-formatted rows below the table ---------------
$ws.Range("D10:D14").NumberFormat = "m/d/yy;@"

# --- Selection / view ------------------------------------------------------
[void]$ws.Range("F11").Select()
